$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the "nome" values in rows 2 and 3, replacing the old placeholder
# name "DAVID AMANCIO...." with "BRUNO SANTA RITA MOREIRA".
$ws.Range("A2").Value = "BRUNO SANTA RITA MOREIRA"
$ws.Range("A3").Value = "THIAGO MARTINS AMORIM"

# Match the final selection / active cell recorded in the saved workbook.
$ws.Range("A7").Select()
